$d = $word.ActiveDocument

$d.Content.Find.Execute("Lucaskepler.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://lucaskepler.netlify.app", 2)
